$wb = $excel.ActiveWorkbook

# --- "data" sheet: add new framework/book columns D:K ---
$wsData = $wb.Worksheets.Item("data")

$cols = @("D","E","F","G","H","I","J","K")
$headers = @("energy_led","low_emit_mat","electric_cars","low_emit_gas","borehole_water","recycle","low_emit_inhale","local_procure")
for ($j = 0; $j -lt $cols.Length; $j++) {
    $wsData.Range($cols[$j] + "1").Value = $headers[$j]
}

$newVals = @(0.33792134495884429, 0.67735511513227553, 0.5, 0.12954974111110007, 0.72609968404906466, 0.22184458107556482, 0.7685331628046782, 0)
for ($i = 2; $i -le 9; $i++) {
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $wsData.Range($cols[$j] + $i).Value = $newVals[$j]
    }
}

# --- "unit_costs" sheet: freeze the volatile RAND() formulas into static values ---
$wsCosts = $wb.Worksheets.Item("unit_costs")

$costCols = @("B","C","D","E","F","G","H","I")
$costVals = @(6170.3416368941998, 4023.5092935694001, 600, 8380.5140125266007, 7867.3384506470002, 1983.2454556056, 8496.1577543293006, 4901.3646667208413)
for ($i = 2; $i -le 9; $i++) {
    for ($j = 0; $j -lt $costCols.Length; $j++) {
        $wsCosts.Range($costCols[$j] + $i).Value = $costVals[$j]
    }
}

# --- Selections / active sheet bookkeeping ---
$wsData.Range("K17").Select()

$wsCosts.Activate()
$wsCosts.Range("B2:I2").Select()
